$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item("Hoja1")

# Update "Periodo Mora" column (E16:E30) from descending (2104..2002)
# to ascending (2002..2104) order. Values are stored as text.
$periodos = @("2002","2003","2004","2005","2006","2007","2008","2009","2010","2011","2012","2101","2102","2103","2104")
for ($i = 0; $i -lt $periodos.Length; $i++) {
    $row = 16 + $i
    $ws.Range("E$row").Value = $periodos[$i]
}

# Swap the two "Valor Mora" date-serial values between the first and
# last data rows (F16 and F30).
$ws.Range("F16").Value = 44944
$ws.Range("F30").Value = 38951
